$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 19:56"

$ws.Range("B4").Value = 3250705
$ws.Range("C4").Value = 30706
$ws.Range("D4").Value = 1437774
$ws.Range("E4").Value = 1676773
$ws.Range("G4").Value = 336
$ws.Range("H4").Value = 136158

$ws.Range("B5").Value = 1768970
$ws.Range("C5").Value = 9867
$ws.Range("E5").Value = 513968
$ws.Range("G5").Value = 152
$ws.Range("H5").Value = 69406

$ws.Range("B6").Value = 821458
$ws.Range("C6").Value = 26616
$ws.Range("D6").Value = 516192
$ws.Range("E6").Value = 283123
$ws.Range("G6").Value = 520
$ws.Range("H6").Value = 22143

$ws.Range("B18").Value = 210965
$ws.Range("C18").Value = 1003
$ws.Range("D18").Value = 191883
$ws.Range("E18").Value = 13759
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 5323

$ws.Range("D19").Value = 184000
$ws.Range("E19").Value = 6206

$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 65018
$ws.Range("C31").Value = 797
$ws.Range("D31").Value = 29577
$ws.Range("E31").Value = 30502
$ws.Range("G31").Value = 39
$ws.Range("H31").Value = 4939

$ws.Range("A32").Value = "Bielorrusia"
$ws.Range("B32").Value = 64604
$ws.Range("C32").Value = 193
$ws.Range("D32").Value = 54254
$ws.Range("E32").Value = 9896
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 454

$ws.Range("B56").Value = 25589
$ws.Range("C56").Value = 24
$ws.Range("E56").Value = 481
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 1744

$ws.Range("B65").Value = 15328
$ws.Range("C65").Value = 249
$ws.Range("D65").Value = 11827
$ws.Range("E65").Value = 3258

$ws.Range("B84").Value = 7120
$ws.Range("C84").Value = 147
$ws.Range("E84").Value = 4566
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 124

$ws.Range("E87").Value = 4436
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 26

$ws.Range("B92").Value = 5704
$ws.Range("C92").Value = 146
$ws.Range("D92").Value = 2785
$ws.Range("E92").Value = 2896
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 23

$ws.Range("B107").Value = 2711
$ws.Range("C107").Value = 9
$ws.Range("E107").Value = 194

$ws.Range("B109").Value = 2617
$ws.Range("C109").Value = 64
$ws.Range("D109").Value = 2238
$ws.Range("E109").Value = 366

$ws.Range("B110").Value = 2451
$ws.Range("C110").Value = 297
$ws.Range("E110").Value = 460

$ws.Range("A111").Value = "Cuba"
$ws.Range("B111").Value = 2413
$ws.Range("C111").Value = 10
$ws.Range("D111").Value = 2249
$ws.Range("E111").Value = 78
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 86

$ws.Range("A112").Value = "Mali"
$ws.Range("B112").Value = 2404
$ws.Range("C112").Value = 34
$ws.Range("D112").Value = 1650
$ws.Range("E112").Value = 633
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 121

$ws.Range("A114").Value = "Malaui"
$ws.Range("B114").Value = 2069
$ws.Range("C114").Value = 127
$ws.Range("D114").Value = 379
$ws.Range("E114").Value = 1659
$ws.Range("G114").Value = 6
$ws.Range("H114").Value = 31

$ws.Range("A115").Value = "Sudan del Sur"
$ws.Range("B115").Value = 2021
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 333
$ws.Range("E115").Value = 1650
$ws.Range("H115").Value = 38

$ws.Range("A116").Value = "Estonia"
$ws.Range("B116").Value = 2013
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 1894
$ws.Range("E116").Value = 50
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 69

$ws.Range("B118").Value = 1886
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = 1859
$ws.Range("E118").Value = 17

$ws.Range("B125").Value = 1591
$ws.Range("C125").Value = 39
$ws.Range("E125").Value = 842
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 19

$ws.Range("B128").Value = 1380
$ws.Range("C128").Value = 24
$ws.Range("D128").Value = 630
$ws.Range("E128").Value = 386
$ws.Range("G128").Value = 3
$ws.Range("H128").Value = 364

$ws.Range("A131").Value = "Suazilandia"
$ws.Range("B131").Value = 1257
$ws.Range("C131").Value = 44
$ws.Range("D131").Value = 633
$ws.Range("E131").Value = 606
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 18

$ws.Range("A132").Value = "Tunez"
$ws.Range("B132").Value = 1240
$ws.Range("C132").Value = 9
$ws.Range("D132").Value = 1067
$ws.Range("E132").Value = 123
$ws.Range("H132").Value = 50

$ws.Range("B139").Value = 1013
$ws.Range("C139").Value = 3
$ws.Range("E139").Value = 155

$ws.Range("B156").Value = 571
$ws.Range("C156").Value = 5
$ws.Range("E156").Value = 96

$ws.Range("A158").Value = "Angola"
$ws.Range("B158").Value = 458
$ws.Range("C158").Value = 62
$ws.Range("D158").Value = 117
$ws.Range("E158").Value = 318
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 23

$ws.Range("A159").Value = "Taiwan"
$ws.Range("B159").Value = 451
$ws.Range("C159").Value = 2
$ws.Range("D159").Value = 438
$ws.Range("E159").Value = 6
$ws.Range("H159").Value = 7
